# Apply attendance_reports sync changes to "Session Analysis Results" sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- L7 / L8: updated summary metric values ---
$ws.Range("L7").Value = 33
$ws.Range("L8").Value = 66

# --- G column: swap "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com" ---
$gRows = @(8,9,10,12,14,15,17,34,35,36,38,40,41,43,60,61,62,64,66,67,69,86,87,88,90,92,93,95,112,113,114,116,118,119,121,138,139,140,142,144,145,147)
foreach ($r in $gRows) {
    $ws.Cells.Item($r, 7).Value = "System, dnasr281@gmail.com"
}

# --- P / Q columns: updated Missing / Pending counts for B1A1..B1C2 ---
$ws.Range("P15").Value = 3
$ws.Range("Q15").Value = 5

$ws.Range("P16").Value = 2
$ws.Range("Q16").Value = 5

$ws.Range("P17").Value = 2
$ws.Range("Q17").Value = 5

$ws.Range("P18").Value = 2
$ws.Range("Q18").Value = 5

$ws.Range("P19").Value = 2
$ws.Range("Q19").Value = 5

$ws.Range("P20").Value = 3
$ws.Range("Q20").Value = 5

# --- Rows that flipped from "Pending" (yellow, style 6) to "Not Recorded" (pink, style 4) ---
# Each of these rows sits directly below a row already styled/labeled "Not Recorded",
# so copy that row's formatting (columns A:I) and then update the status text.
$flippedRows = @(22,48,74,100,126,152)
foreach ($r in $flippedRows) {
    $srcRow = $r - 1
    $ws.Range("A" + $srcRow + ":I" + $srcRow).Copy() | Out-Null
    $ws.Range("A" + $r + ":I" + $r).PasteSpecial(-4122) | Out-Null
    $ws.Range("I" + $r).Value = "Not Recorded"
}
